# Fixed #97 and #98
#
# #97: split the old "inline-anchors-in-paragraph" / "inline-anchors-in-tagged-paragraph"
#      rows into differently-named/tagged rows, and re-home "rule_with_newlines" under
#      the "paragraph" norm namespace, right after the other paragraph rows.
# #98: fix the unordered-list / description-list test fixture text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Normative Rules")

# Row 2 ("inline") is unchanged.

# Row 3 used to be "inline-anchors-in-paragraph"; it is now the (shifted-up)
# "paragraph-with-a-really-wide-rule-name" row, whose tag text also changed.
$ws.Cells.Item(3, 1).Value = "my-chapter_name"
$ws.Cells.Item(3, 2).Value = "paragraph-with-a-really-wide-rule-name"
$ws.Cells.Item(3, 3).Value = "Here's a description.`nIt's got 2 lines.`nParagraph without inline anchors"
$ws.Cells.Item(3, 4).Value = 'Description, ["norm:paragraph:no-inline-anchors"]'

# Row 4: new "entire paragraph" anchors case (replaces old inline-anchors-in-paragraph row).
$ws.Cells.Item(4, 1).Value = "my-chapter_name"
$ws.Cells.Item(4, 2).Value = "inline-anchors-in-paragraph-entire"
$ws.Cells.Item(4, 3).Value = "Paragraph with inline anchor and something."
$ws.Cells.Item(4, 4).Value = '["norm:paragraph:inline-anchors:entire"]'

# Row 5: new "inline anchor" case (replaces old inline-anchors-in-tagged-paragraph row).
$ws.Cells.Item(5, 1).Value = "my-chapter_name"
$ws.Cells.Item(5, 2).Value = "inline-anchors-in-paragraph-inline-anchor"
$ws.Cells.Item(5, 3).Value = "inline anchor"
$ws.Cells.Item(5, 4).Value = '["norm:paragraph:inline-anchors:inline-anchor"]'

# Row 6: "rule_with_newlines" moved up here (was row 12) and re-tagged under "paragraph".
$ws.Cells.Item(6, 1).Value = "my-chapter_name"
$ws.Cells.Item(6, 2).Value = "rule_with_newlines"
$ws.Cells.Item(6, 3).Value = "Here&#8217;s the first line. Here&#8217;s the second line."
$ws.Cells.Item(6, 4).Value = '["norm:paragraph:tag_with_newlines"]'

# Row 7: "table1" shifted down one row (from old row 6); its formula cell (C7) is left
# untouched below so the stored formula text is preserved as-is.
$ws.Cells.Item(7, 1).Value = "my-chapter_name"
$ws.Cells.Item(7, 2).Value = "table1"
$ws.Cells.Item(7, 4).Value = '["norm:table:anchors-in-cells:entire-table"]'
$ws.Cells.Item(7, 3).Formula = "===`n cell with anchor`ncell without anchor`n==="

# Row 8: "table2" shifted down one row (from old row 7).
$ws.Cells.Item(8, 1).Value = "my-chapter_name"
$ws.Cells.Item(8, 2).Value = "table2"
$ws.Cells.Item(8, 4).Value = '["norm:table:no-anchors-in-cells:entire-table"]'
$ws.Cells.Item(8, 3).Formula = "=Header 1|Header 2`n==`nCell in column 1, row 1|Cell in column 2, row 1`nCell in column 1, row 2|Cell in column 2, row 2`n==="

# Row 9: "unordered1" shifted down one row (from old row 8); fixture text updated (#98).
# Old row 9 (note_with_2_tags) had E9/F9 populated; clear them since this row no longer
# needs a 5th/6th column.
$ws.Cells.Item(9, 1).Value = "my-chapter_name"
$ws.Cells.Item(9, 2).Value = "unordered1"
$ws.Cells.Item(9, 3).Value = "Item 1`n Item 2`nItem 3"
$ws.Cells.Item(9, 4).Value = '["norm:unordered-list:anchors-in-items:entire-list"]'
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(9, 6).ClearContents()

# Row 10: "note_with_2_tags" shifted down one row (from old row 9); content unchanged.
$ws.Cells.Item(10, 1).Value = "my-chapter_name"
$ws.Cells.Item(10, 2).Value = "note_with_2_tags"
$ws.Cells.Item(10, 3).Value = "One line description`nParagraph 1`nParagraph 3"
$ws.Cells.Item(10, 4).Value = 'Description, ["norm:admonition:anchors-in-notes:note1", "norm:admonition:anchors-in-notes:note3"]'
$ws.Cells.Item(10, 5).Value = "parameter"
$ws.Cells.Item(10, 6).Value = "MY_PARAMETER"

# Row 11: "desc1" shifted down one row (from old row 10); fixture text updated (#98).
$ws.Cells.Item(11, 1).Value = "my-chapter_name"
$ws.Cells.Item(11, 2).Value = "desc1"
$ws.Cells.Item(11, 3).Value = "Item 1`nItem 3"
$ws.Cells.Item(11, 4).Value = '["norm:description-list:anchors-in-items:item1", "norm:description-list:anchors-in-items:item3"]'

# Row 12: "desc2" (same row position as before); fixture text updated (#98).
$ws.Cells.Item(12, 1).Value = "my-chapter_name"
$ws.Cells.Item(12, 2).Value = "desc2"
$ws.Cells.Item(12, 3).Value = "Item 1`nItem 3"
$ws.Cells.Item(12, 4).Value = '["norm:description-list:anchors-in-items:item1", "norm:description-list:anchors-in-items:item3"]'
